$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.942.92'
$ws.Range('E2').Value = '  -1.79%  '
$ws.Range('D3').Value = '1.635.28'
$ws.Range('E3').Value = '  -1.97%  '
$ws.Range('E4').Value = '  +0.58%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '215.69'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -1.59%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '1.013'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +0.59%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.5001'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -3.45%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.2566'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  -0.39%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.06419'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  -0.67%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '19.47'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -2.61%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.07747'
$cell.Style = 'Normal'
$ws.Range('D12').Value = '1.641.73'
$ws.Range('E12').Value = '  -1.62%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '4.252'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  -2.25%  '
$ws.Range('D14').Value = '1.862.36'
$ws.Range('E14').Value = '  -1.87%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '0.5441'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -2.12%  '
$ws.Range('D16').Value = '0.0₅7920'
$ws.Range('E16').Value = '  -1.80%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '63.52'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('D18').Value = '25.965.41'
$ws.Range('E18').Value = '  -1.78%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '1.013'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +0.58%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '203.18'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -3.51%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '4.298'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -2.92%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '9.997'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -1.32%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '5.952'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +0.84%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '1.013'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +0.64%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '1.975'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +14.31%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '140.94'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -2.80%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '0.1148'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -1.84%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '15.73'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -0.74%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '6.797'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -3.14%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '0.05061'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -3.78%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '1.243'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -1.67%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '3.260'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -3.69%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '3.196'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -0.94%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '1.545'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -2.07%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '2.348'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -1.15%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.8916'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -4.13%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '2.607'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  -5.61%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.5641'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -1.72%  '
$ws.Range('D39').Value = '1.133.21'
$ws.Range('E39').Value = '  -1.73%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.01560'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -2.97%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '2.586'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +0.22%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '1.013'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +0.67%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '5.639'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -0.33%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '0.8179'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -3.70%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '99.75'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -0.53%  '
$ws.Range('D46').Value = '1.773.75'
$ws.Range('E46').Value = '  -1.85%  '
$ws.Range('D47').Value = '0.0₈112'
$ws.Range('E47').Value = '  +0.56%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '0.4538'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +1.00%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '1.013'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +0.53%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '54.75'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -2.28%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '0.05029'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -1.68%  '
